# Applies the "added origin_destination folders in input and test/input,
# also added in config files" edit to the `input` sheet of user_input_tests.xlsx.
#
# Concretely (per the canonical-XML diff):
#   - input!D2  "Network based on OSM online" -> "Network based on shapefile"
#   - input!G2  "zuidholland_4326"            -> "part_of_DR_roads"
#   - input!H2  "osmid"                       -> "fid"
#   - the `input` sheet becomes the active/selected sheet/tab, with the
#     cursor resting on D14 (was D2)
#   - the `explanation` sheet loses its "tabSelected" flag (only one sheet
#     can be the active tab - a side effect of activating `input` instead)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("input")

# --- cell value updates (row 2) -----------------------------------------
$ws.Range("D2").Value = "Network based on shapefile"
$ws.Range("G2").Value = "part_of_DR_roads"
$ws.Range("H2").Value = "fid"

# --- selection / active-sheet updates -----------------------------------
$ws.Activate()
$ws.Range("D14").Select()
